# Attendance.xlsx - "Added new entries to Attendance"
#
# Adds:
#   - C3: "Aanwezig - " (a missed entry for Liam on the existing 2nd tracked day)
#   - Four brand-new attendance rows (4-7) for the weeks of 19/20 and 26/27 May 2025
#
# Shared-string cache ordering matters for byte-identical sharedStrings.xml, so
# cells are written in natural row-major / left-to-right order (row 3, then
# rows 4-7, columns A..G) matching the order values first appear.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in the previously-empty "Liam" cell -----------------------
$ws.Range("C3").Value = "Aanwezig - "

# --- Row 4 (2025-05-19 / Monday) -------------------------------------------
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat
$ws.Range("A4").Value = 45796
$ws.Range("B4").Value = "Aanwezig - Class Diagram"
$ws.Range("C4").Value = "Aanwezig - "
$ws.Range("D4").Value = "Aanwezig - Niew lijst toevoegen"
$ws.Range("E4").Value = "Afwezig"
$ws.Range("F4").Value = "Ziek"
$ws.Range("G4").Value = "Aanwezig - Onduidelijk"

# --- Row 5 (2025-05-20 / Tuesday) -------------------------------------------
$ws.Range("A5").NumberFormat = $ws.Range("A3").NumberFormat
$ws.Range("A5").Value = 45797
$ws.Range("B5").Value = "Aanwezig - Class Diagram"
$ws.Range("C5").Value = "Aanwezig - "
$ws.Range("D5").Value = "Aanwezig - Niew lijst toevoegen"
$ws.Range("E5").Value = "Afwezig"
$ws.Range("F5").Value = "Ziek"
$ws.Range("G5").Value = "Aanwezig - Onduidelijk"

# --- Row 6 (2025-05-26 / Monday) -------------------------------------------
$ws.Range("A6").NumberFormat = $ws.Range("A3").NumberFormat
$ws.Range("A6").Value = 45803
$ws.Range("B6").Value = "Aanwezig - Class Diagram"
$ws.Range("C6").Value = "Aanwezig - "
$ws.Range("D6").Value = "Aanwezig - Product aan lijst toevoegen"
$ws.Range("E6").Value = "Afwezig"
$ws.Range("F6").Value = "Aanwezig - Product aan lijst toevoegen"
$ws.Range("G6").Value = "Aanwezig - Onduidelijk"

# --- Row 7 (2025-05-27 / Tuesday) -------------------------------------------
$ws.Range("A7").NumberFormat = $ws.Range("A3").NumberFormat
$ws.Range("A7").Value = 45804
$ws.Range("B7").Value = "Aanwezig - Class Diagram"
$ws.Range("C7").Value = "Aanwezig - "
$ws.Range("D7").Value = "Aanwezig - UI van pagina maken"
$ws.Range("E7").Value = "Afwezig"
$ws.Range("F7").Value = "Aanwezig - Product aan lijst toevoegen"
$ws.Range("G7").Value = "Afwezig"

# --- Match the author's final cursor position -------------------------------
$ws.Range("D12").Select()
